$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5: 45183 -> 45184
$ws.Range("C2").Value = 45184
$ws.Range("C3").Value = 45184
$ws.Range("C4").Value = 45184
$ws.Range("C5").Value = 45184
